$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-13 Saturday", "2024-04-14 Sunday"),
    @("426×2=", "425×3="),
    @("386×9=", "187×7="),
    @("309×4=", "127×2="),
    @("615×7=", "833×9="),
    @("525×5=", "671×3="),
    @("665×9=", "110×9="),
    @("389×9=", "815×5="),
    @("345×8=", "507×6="),
    @("420×2=", "788×5="),
    @("171×7=", "747×8="),
    @("716×7=", "754×7="),
    @("474×2=", "641×9="),
    @("343×6=", "265×4="),
    @("614×2=", "627×3="),
    @("824×4=", "118×2="),
    @("150×2=", "510×9="),
    @("643×9=", "405×8="),
    @("849×2=", "690×4="),
    @("161×2=", "739×3="),
    @("675×4=", "171×4="),
    @("417×3=", "979×8="),
    @("543×5=", "504×9="),
    @("765×7=", "347×3="),
    @("447×7=", "797×3="),
    @("991×8=", "115×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
